# Natmi following Dr Hou advice
#
# Recomputes the L1cam-Itgav sending/target-cluster table to include a new
# "FAPs" cluster (in addition to the existing ECs, M2 and sCs clusters),
# expanding the cross-joined Sending cluster x Target cluster rows from
# 3x4=12 to 4x4=16 and refreshing all of the associated NATMI statistics.
#
# Columns A:T are, in order:
#   A Sending cluster            K Receptor-expressing cells
#   B Ligand symbol              L Receptor detection rate
#   C Receptor symbol            M Receptor average expression value
#   D Target cluster             N Receptor total expression value
#   E Ligand-expressing cells    O Receptor derived specificity (avg)
#   F Ligand detection rate      P Receptor derived specificity (total)
#   G Ligand average expr value Q Edge average expression weight
#   H Ligand total expr value   R Edge total expression weight
#   I Ligand derived spec (avg) S Edge average expression derived spec.
#   J Ligand derived spec (tot) T Edge total expression derived spec.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "L1cam", "Itgav", "ECs", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 159.778035800779, 1438.002322207011, 0.0535444385467083, 0.05354443854670831),
    @("ECs", "L1cam", "Itgav", "FAPs", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 650.3628190125182, 5853.265371112664, 0.2179480541311672, 0.2179480541311672),
    @("ECs", "L1cam", "Itgav", "M2", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 756.9377237573433, 6812.43951381609, 0.2536631848694483, 0.2536631848694483),
    @("ECs", "L1cam", "Itgav", "sCs", 3, 1, 11.88712433333333, 35.661373, 0.5967229292030898, 0.5967229292030898, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 213.5585918459742, 1922.027326613768, 0.07156725165576591, 0.07156725165576593),
    @("FAPs", "L1cam", "Itgav", "ECs", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 5.088725550287, 45.798529952583, 0.001705321705469847, 0.001705321705469847),
    @("FAPs", "L1cam", "Itgav", "FAPs", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 20.71322179847166, 186.418996186245, 0.006941366040295243, 0.006941366040295244),
    @("FAPs", "L1cam", "Itgav", "M2", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 24.10749584919667, 216.96746264277, 0.008078847155323728, 0.00807884715532373),
    @("FAPs", "L1cam", "Itgav", "sCs", 2, 0.6666666666666666, 0.3785896666666667, 1.135769, 0.01900486009268527, 0.01900486009268527, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 6.801567295300445, 61.214105657704, 0.002279325191596454, 0.002279325191596454),
    @("M2", "L1cam", "Itgav", "ECs", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 55.421013458262, 498.789121124358, 0.01857255932856908, 0.01857255932856909),
    @("M2", "L1cam", "Itgav", "FAPs", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 225.58649168893, 2030.27842520037, 0.07559801308526626, 0.07559801308526629),
    @("M2", "L1cam", "Itgav", "M2", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 262.55333259778, 2362.97999338002, 0.08798625363027997, 0.08798625363027998),
    @("M2", "L1cam", "Itgav", "sCs", 3, 1, 4.123197999999999, 12.369594, 0.2069808238940481, 0.2069808238940482, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 74.07547309932266, 666.679257893904, 0.02482399784993281, 0.02482399784993282),
    @("sCs", "L1cam", "Itgav", "ECs", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 13.441269, 40.323807, 0.08973082133481231, 0.08973082133481232, 47.471394448939, 427.242550040451, 0.01590850175406507, 0.01590850175406507),
    @("sCs", "L1cam", "Itgav", "FAPs", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 54.711535, 164.134605, 0.3652416280068742, 0.3652416280068742, 193.2282479349183, 1739.054231414265, 0.06475419475014542, 0.06475419475014542),
    @("sCs", "L1cam", "Itgav", "M2", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 63.67711, 191.03133, 0.4250937452800914, 0.4250937452800915, 224.8925459477433, 2024.03291352969, 0.07536545962503943, 0.07536545962503945),
    @("sCs", "L1cam", "Itgav", "sCs", 3, 1, 3.531764333333333, 10.595293, 0.1772913868101768, 0.1772913868101768, 3, 1, 17.96553866666667, 53.896616, 0.119933805378222, 0.119933805378222, 63.45004869205422, 571.050438228488, 0.02126323068092681, 0.02126323068092682)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}
